$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting existing rows 4-10 down to 5-11.
$ws.Rows(4).Insert()

# Fill the newly inserted row 4 with the new weekly data point.
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 44804
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 100112012
$ws.Range("G4").Value = "Espinaca"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 7000
$ws.Range("L4").Value = 7500
$ws.Range("M4").Value = 7250
$ws.Range("N4").Value = "$/cuna 10 kilos"
$ws.Range("O4").Value = "Provincia de Diguillín"
$ws.Range("P4").Value = 725
$ws.Range("Q4").Value = 10
$ws.Range("R4").Value = "Hortaliza"

# Safety net: ensure the Origen value for the (shifted) row 9 is correct.
$ws.Range("O9").Value = "Región Metropolitana"
